# Natmi following Dr Hou advice
# Recomputes the Nppc-Npr3 LR-pair table: adds a new "ECs" sending/target
# cluster, corrects the FAPs/sCs row pairing, and refreshes every derived
# expression/specificity metric (rows 2-7, columns E-T).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Nppc"
$ws.Cells.Item(2, 3).Value = "Npr3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.345980666666667
$ws.Cells.Item(2, 8).Value = 4.037942
$ws.Cells.Item(2, 9).Value = 0.4408299556445331
$ws.Cells.Item(2, 10).Value = 0.4408299556445331
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.2990286666666667
$ws.Cells.Item(2, 14).Value = 0.897086
$ws.Cells.Item(2, 15).Value = 0.2003670139510866
$ws.Cells.Item(2, 16).Value = 0.2003670139510866
$ws.Cells.Item(2, 17).Value = 0.4024868041124445
$ws.Cells.Item(2, 18).Value = 3.622381237012001
$ws.Cells.Item(2, 19).Value = 0.08832778187268504
$ws.Cells.Item(2, 20).Value = 0.08832778187268506
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Nppc"
$ws.Cells.Item(3, 3).Value = "Npr3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.345980666666667
$ws.Cells.Item(3, 8).Value = 4.037942
$ws.Cells.Item(3, 9).Value = 0.4408299556445331
$ws.Cells.Item(3, 10).Value = 0.4408299556445331
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.007656
$ws.Cells.Item(3, 14).Value = 3.022968
$ws.Cells.Item(3, 15).Value = 0.675189526343838
$ws.Cells.Item(3, 16).Value = 0.675189526343838
$ws.Cells.Item(3, 17).Value = 1.356285494650667
$ws.Cells.Item(3, 18).Value = 12.206569451856
$ws.Cells.Item(3, 19).Value = 0.2976437689498074
$ws.Cells.Item(3, 20).Value = 0.2976437689498074
$ws.Cells.Item(4, 1).Value = "sCs"
$ws.Cells.Item(4, 2).Value = "Nppc"
$ws.Cells.Item(4, 3).Value = "Npr3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.345980666666667
$ws.Cells.Item(4, 8).Value = 4.037942
$ws.Cells.Item(4, 9).Value = 0.4408299556445331
$ws.Cells.Item(4, 10).Value = 0.4408299556445331
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.18572
$ws.Cells.Item(4, 14).Value = 0.55716
$ws.Cells.Item(4, 15).Value = 0.1244434597050755
$ws.Cells.Item(4, 16).Value = 0.1244434597050755
$ws.Cells.Item(4, 17).Value = 0.2499755294133333
$ws.Cells.Item(4, 18).Value = 2.24977976472
$ws.Cells.Item(4, 19).Value = 0.05485840482204069
$ws.Cells.Item(4, 20).Value = 0.0548584048220407
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Nppc"
$ws.Cells.Item(5, 3).Value = "Npr3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.707307
$ws.Cells.Item(5, 8).Value = 5.121921
$ws.Cells.Item(5, 9).Value = 0.5591700443554668
$ws.Cells.Item(5, 10).Value = 0.5591700443554669
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.2990286666666667
$ws.Cells.Item(5, 14).Value = 0.897086
$ws.Cells.Item(5, 15).Value = 0.2003670139510866
$ws.Cells.Item(5, 16).Value = 0.2003670139510866
$ws.Cells.Item(5, 17).Value = 0.5105337358006666
$ws.Cells.Item(5, 18).Value = 4.594803622206
$ws.Cells.Item(5, 19).Value = 0.1120392320784015
$ws.Cells.Item(5, 20).Value = 0.1120392320784016
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Nppc"
$ws.Cells.Item(6, 3).Value = "Npr3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.707307
$ws.Cells.Item(6, 8).Value = 5.121921
$ws.Cells.Item(6, 9).Value = 0.5591700443554668
$ws.Cells.Item(6, 10).Value = 0.5591700443554669
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.007656
$ws.Cells.Item(6, 14).Value = 3.022968
$ws.Cells.Item(6, 15).Value = 0.675189526343838
$ws.Cells.Item(6, 16).Value = 0.675189526343838
$ws.Cells.Item(6, 17).Value = 1.720378142392
$ws.Cells.Item(6, 18).Value = 15.483403281528
$ws.Cells.Item(6, 19).Value = 0.3775457573940305
$ws.Cells.Item(6, 20).Value = 0.3775457573940306
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Nppc"
$ws.Cells.Item(7, 3).Value = "Npr3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.707307
$ws.Cells.Item(7, 8).Value = 5.121921
$ws.Cells.Item(7, 9).Value = 0.5591700443554668
$ws.Cells.Item(7, 10).Value = 0.5591700443554669
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.18572
$ws.Cells.Item(7, 14).Value = 0.55716
$ws.Cells.Item(7, 15).Value = 0.1244434597050755
$ws.Cells.Item(7, 16).Value = 0.1244434597050755
$ws.Cells.Item(7, 17).Value = 0.31708105604
$ws.Cells.Item(7, 18).Value = 2.85372950436
$ws.Cells.Item(7, 19).Value = 0.06958505488303483
$ws.Cells.Item(7, 20).Value = 0.06958505488303485
